# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Refresh the worker/period detail table (rows 16-30). Table layout:
#   B = Tipo Doc Trabajador (unchanged, stays "CC")
#   C = N Doc Trabajador
#   D = Nombre Trabajador
#   E = Periodo Mora
#   F = Valor Mora (unchanged, stays 40000)
#   G = Salario Basico

$rows = @(
    @{ Row = 16; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2309"; Salario = 1000000 },
    @{ Row = 17; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2209"; Salario = 1000000 },
    @{ Row = 18; Doc = "30879583";  Nombre = "EIMY DE LOS MILAGROS BARRAZA PINTO";  Periodo = "2209"; Salario = 1000000 },
    @{ Row = 19; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2210"; Salario = 1000000 },
    @{ Row = 20; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2211"; Salario = 1000000 },
    @{ Row = 21; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2212"; Salario = 1000000 },
    @{ Row = 22; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2301"; Salario = 1000000 },
    @{ Row = 23; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2302"; Salario = 1000000 },
    @{ Row = 24; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2303"; Salario = 1000000 },
    @{ Row = 25; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2304"; Salario = 1000000 },
    @{ Row = 26; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2305"; Salario = 1000000 },
    @{ Row = 27; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2306"; Salario = 1000000 },
    @{ Row = 28; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2307"; Salario = 1000000 },
    @{ Row = 29; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2308"; Salario = 1000000 },
    @{ Row = 30; Doc = "52706194";  Nombre = "YOLADIS PUELLO LLERENA";              Periodo = "2208"; Salario = 1000000 }
)

foreach ($r in $rows) {
    $ws.Range("C" + $r.Row).Value = $r.Doc
    $ws.Range("D" + $r.Row).Value = $r.Nombre
    $ws.Range("E" + $r.Row).Value = $r.Periodo
    $ws.Range("G" + $r.Row).Value = $r.Salario
}
